$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift AgTests (H) / AgPosit (I) values up by one row for rows 220-284,
# and set the new values for row 285 (previously the last row).
$ws.Cells.Item(220,8).Value = 0
$ws.Cells.Item(220,9).Value = 0
$ws.Cells.Item(221,8).Value = 525
$ws.Cells.Item(221,9).Value = 40
$ws.Cells.Item(222,8).Value = 1155
$ws.Cells.Item(222,9).Value = 84
$ws.Cells.Item(223,8).Value = 929
$ws.Cells.Item(223,9).Value = 31
$ws.Cells.Item(224,8).Value = 372
$ws.Cells.Item(224,9).Value = 36
$ws.Cells.Item(225,8).Value = 320
$ws.Cells.Item(225,9).Value = 27
$ws.Cells.Item(226,8).Value = 0
$ws.Cells.Item(226,9).Value = 0
$ws.Cells.Item(227,8).Value = 92
$ws.Cells.Item(227,9).Value = 0
$ws.Cells.Item(228,8).Value = 792
$ws.Cells.Item(228,9).Value = 69
$ws.Cells.Item(229,8).Value = 613
$ws.Cells.Item(229,9).Value = 38
$ws.Cells.Item(230,8).Value = 954
$ws.Cells.Item(230,9).Value = 48
$ws.Cells.Item(231,8).Value = 2136
$ws.Cells.Item(231,9).Value = 141
$ws.Cells.Item(232,8).Value = 2161
$ws.Cells.Item(232,9).Value = 119
$ws.Cells.Item(233,8).Value = 1194
$ws.Cells.Item(233,9).Value = 65
$ws.Cells.Item(234,8).Value = 729
$ws.Cells.Item(234,9).Value = 40
$ws.Cells.Item(235,8).Value = 2161
$ws.Cells.Item(235,9).Value = 231
$ws.Cells.Item(236,8).Value = 3036
$ws.Cells.Item(236,9).Value = 234
$ws.Cells.Item(237,8).Value = 2415
$ws.Cells.Item(237,9).Value = 225
$ws.Cells.Item(238,8).Value = 6149
$ws.Cells.Item(238,9).Value = 239
$ws.Cells.Item(239,8).Value = 41133
$ws.Cells.Item(239,9).Value = 471
$ws.Cells.Item(240,8).Value = 89674
$ws.Cells.Item(240,9).Value = 975
$ws.Cells.Item(241,8).Value = 30176
$ws.Cells.Item(241,9).Value = 376
$ws.Cells.Item(242,8).Value = 26701
$ws.Cells.Item(242,9).Value = 1331
$ws.Cells.Item(243,8).Value = 5279
$ws.Cells.Item(243,9).Value = 98
$ws.Cells.Item(244,8).Value = 3592
$ws.Cells.Item(244,9).Value = 82
$ws.Cells.Item(245,8).Value = 2115
$ws.Cells.Item(245,9).Value = 140
$ws.Cells.Item(246,8).Value = 34523
$ws.Cells.Item(246,9).Value = 327
$ws.Cells.Item(247,8).Value = 45190
$ws.Cells.Item(247,9).Value = 439
$ws.Cells.Item(248,8).Value = 12172
$ws.Cells.Item(248,9).Value = 143
$ws.Cells.Item(249,8).Value = 17651
$ws.Cells.Item(249,9).Value = 1154
$ws.Cells.Item(250,8).Value = 6492
$ws.Cells.Item(250,9).Value = 140
$ws.Cells.Item(251,8).Value = 4620
$ws.Cells.Item(251,9).Value = 101
$ws.Cells.Item(252,8).Value = 5606
$ws.Cells.Item(252,9).Value = 173
$ws.Cells.Item(253,8).Value = 6144
$ws.Cells.Item(253,9).Value = 158
$ws.Cells.Item(254,8).Value = 5452
$ws.Cells.Item(254,9).Value = 99
$ws.Cells.Item(255,8).Value = 1242
$ws.Cells.Item(255,9).Value = 33
$ws.Cells.Item(256,8).Value = 5258
$ws.Cells.Item(256,9).Value = 275
$ws.Cells.Item(257,8).Value = 3874
$ws.Cells.Item(257,9).Value = 220
$ws.Cells.Item(258,8).Value = 6114
$ws.Cells.Item(258,9).Value = 423
$ws.Cells.Item(259,8).Value = 11540
$ws.Cells.Item(259,9).Value = 750
$ws.Cells.Item(260,8).Value = 17458
$ws.Cells.Item(260,9).Value = 597
$ws.Cells.Item(261,8).Value = 9226
$ws.Cells.Item(261,9).Value = 317
$ws.Cells.Item(262,8).Value = 1701
$ws.Cells.Item(262,9).Value = 57
$ws.Cells.Item(263,8).Value = 42064
$ws.Cells.Item(263,9).Value = 860
$ws.Cells.Item(264,8).Value = 17221
$ws.Cells.Item(264,9).Value = 854
$ws.Cells.Item(265,8).Value = 13146
$ws.Cells.Item(265,9).Value = 676
$ws.Cells.Item(266,8).Value = 13970
$ws.Cells.Item(266,9).Value = 809
$ws.Cells.Item(267,8).Value = 15183
$ws.Cells.Item(267,9).Value = 716
$ws.Cells.Item(268,8).Value = 9461
$ws.Cells.Item(268,9).Value = 413
$ws.Cells.Item(269,8).Value = 2672
$ws.Cells.Item(269,9).Value = 177
$ws.Cells.Item(270,8).Value = 41562
$ws.Cells.Item(270,9).Value = 1537
$ws.Cells.Item(271,8).Value = 30827
$ws.Cells.Item(271,9).Value = 1693
$ws.Cells.Item(272,8).Value = 25977
$ws.Cells.Item(272,9).Value = 1301
$ws.Cells.Item(273,8).Value = 27350
$ws.Cells.Item(273,9).Value = 1260
$ws.Cells.Item(274,8).Value = 27135
$ws.Cells.Item(274,9).Value = 1176
$ws.Cells.Item(275,8).Value = 12360
$ws.Cells.Item(275,9).Value = 400
$ws.Cells.Item(276,8).Value = 3069
$ws.Cells.Item(276,9).Value = 120
$ws.Cells.Item(277,8).Value = 29022
$ws.Cells.Item(277,9).Value = 1977
$ws.Cells.Item(278,8).Value = 42220
$ws.Cells.Item(278,9).Value = 2946
$ws.Cells.Item(279,8).Value = 32670
$ws.Cells.Item(279,9).Value = 2130
$ws.Cells.Item(280,8).Value = 41769
$ws.Cells.Item(280,9).Value = 2948
$ws.Cells.Item(281,8).Value = 42696
$ws.Cells.Item(281,9).Value = 2581
$ws.Cells.Item(282,8).Value = 16398
$ws.Cells.Item(282,9).Value = 951
$ws.Cells.Item(283,8).Value = 1509
$ws.Cells.Item(283,9).Value = 118
$ws.Cells.Item(284,8).Value = 33530
$ws.Cells.Item(284,9).Value = 2696
$ws.Cells.Item(285,8).Value = 19328
$ws.Cells.Item(285,9).Value = 1635

# Append the new row 286 with data for 2020-12-15 (no AgTests/AgPosit yet).
$ws.Cells.Item(286,1).Value = 44180
$ws.Cells.Item(286,1).NumberFormat = "yyyy-mm-dd"
$ws.Cells.Item(286,2).Value = 139088
$ws.Cells.Item(286,3).Value = 101584
$ws.Cells.Item(286,4).Value = 36195
$ws.Cells.Item(286,5).Value = 16717
$ws.Cells.Item(286,6).Value = 3565
$ws.Cells.Item(286,7).Value = 1309

